# The records in rows 10-28 of the "Artfynd" sheet were re-matched to
# their correct GPS/time/observer data. For several pairs of rows the
# fix shows up as the two rows trading the values of a specific set of
# columns (the columns that actually differ between the two records),
# while the columns that already agree between the pair are left alone.
#
# Swapping only the columns that differ (instead of the whole row) avoids
# rewriting date-looking text cells (column Y/AA hold "2025-08-14" in
# every one of these rows) through Range.Value, which would otherwise get
# auto-converted from text into a real Excel date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cells($row1, $row2, [string[]]$cols) {
    foreach ($col in $cols) {
        $cellA = $ws.Range("$col$row1")
        $cellB = $ws.Range("$col$row2")

        $valA = $cellA.Value()
        $valB = $cellB.Value()

        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

Swap-Cells 10 11 @("A","I","Q","R","X","Z","AB","AX")
Swap-Cells 13 15 @("A","I","Q","R","X","Z","AB","AX")
Swap-Cells 20 21 @("A","I","J","Q","R","X","Z","AB","AX")
Swap-Cells 23 26 @("A","B","E","F","G","H","I","J","Q","R","X","Z","AB","AC")
Swap-Cells 24 25 @("A","J","Q","R","X","Z","AB","AX")
Swap-Cells 27 28 @("A","B","D","E","F","G","H","J","Q","R","X","Z","AB")
